# Weekly "Fruta / hortaliza" update: a new daily price record is inserted
# for "Achicoria" at row 79 (Vega Central Mapocho de Santiago), pushing the
# existing rows 79-88 down to 80-89. The new row carries the same static
# descriptive fields as its neighbours, but its own date / price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 79, shifting rows 79:88 down to 80:89 (and
# bringing the dimension from A1:R88 to A1:R89).
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(79, 1).Value = 9
$ws.Cells.Item(79, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(79, 3).Value = "Metropolitana"
$ws.Cells.Item(79, 4).Value = 45212
$ws.Cells.Item(79, 5).Value = 13
$ws.Cells.Item(79, 6).Value = 100112010
$ws.Cells.Item(79, 7).Value = "Achicoria"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 90
$ws.Cells.Item(79, 11).Value = 7000
$ws.Cells.Item(79, 12).Value = 7000
$ws.Cells.Item(79, 13).Value = 7000
$ws.Cells.Item(79, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(79, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(79, 16).Value = 438
$ws.Cells.Item(79, 17).Value = 16
$ws.Cells.Item(79, 18).Value = "Hortaliza"
